{"js": "// Remove a set of resume bullet lines (each bullet is a run of text followed\n// by a manual line break run, all inside one larger <w:p> \"story\" paragraph).\n// We search for each bullet's exact text PLUS the manual line break char\n// (\\u000b, how Office.js represents a <w:br/> in range.text) so that a single\n// `.delete()` removes both the text run and its trailing break run, exactly\n// like the target diff (no leftover empty runs/paragraphs).\n\nconst body = context.document.body;\n\nconst bulletsToRemove = [\n  \"\\u2022 Optimized, commercialized, and launched a two\\u00ad-component, polyurethane clearcoat system which delivered $800M revenue growth in 2012.\",\n  \"\\u2022 Developed new rheology test method with optimized shear profile, reducing error in measurement by 50%, for non\\u00ad-Newtonian solventborne paint systems.\",\n  \"\\u2022 Technical lead for innovative spray process and paint technology conversion program which minimized assembly line downtime by 50% over conventional conversion.\",\n  \"\\u2022 Provided support to Arlington Assembly, GM's most profitable manufacturing site, via new color development, formulation adjustments, and troubleshooting line issues (2011-2013). Currently provide support to Bowling Green Assembly, home of the Corvette.\",\n  \"\\u2022 Align formulas, manufacturing procedures, and product design specifications for manufacturing scale\\u00ad-up.\",\n  \"\\u2022 Gravimetric measurement of binary vapor\\u00ad-liquid equilibrium curves of ionic liquids with components in flue gas (e.g. CO2, CH4, H2O) as well as N2O.\",\n  \"\\u2022 Analyzed and calculated hysteresis, infinite dilution activity coefficients, Henry's Law constants, and deconvolution of physical CO2 solubility from chemical CO2 reaction in amine\\u00ad-functionalized ionic liquid systems.\",\n  \"\\u2022 Supervised design and construction of ionic liquid absorber/ stripper unit.\",\n  \"\\u2022 Proficient with both high and low pressure systems.\",\n  \"\\u2022 Facilitated formulation of controlled\\u00ad-release coating.\",\n  \"\\u2022 Authored process flow diagrams. Collaborated on design of specialized fluidized bed dryer. Designated process instrumentation for final scale\\u00adup.\",\n  \"\\u2022 Supervised 3\\u00ad-4 production workers per shift in pilot plant operations.\",\n  \"\\u2022 Designed and formulated a novel filtration system capable of 100% toxin removal from contaminated liquids.\",\n  \"\\u2022 Bottled beer (KY Ale, KY Light, KY Bourbon Barrel Ale) at company microbrewery.\",\n];\n\nfor (const bulletText of bulletsToRemove) {\n  // Include the trailing manual-line-break char so the whole \"line\" (text +\n  // its <w:br/>) is matched and removed together.\n  const searchTarget = bulletText + \"\\u000b\";\n  const results = body.search(searchTarget, {\n    matchCase: true,\n    matchWildcards: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Bullet text not found: \" + bulletText.substring(0, 50));\n  }\n\n  results.items[0].delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each bullet below is one \"line\" inside a larger paragraph: the bullet's\n# text run immediately followed by a manual line-break run (<w:br/>, which\n# Word's Range model exposes as character code 11 / vertical tab). We search\n# for the bullet text and append that break character to the search string so\n# a single Find + Delete removes the text run AND its trailing break run\n# together, leaving no orphaned empty runs (matching the target diff).\n$bulletsToRemove = @(\n    ([char]0x2022 + ' Optimized, commercialized, and launched a two' + [char]0xAD + '-component, polyurethane clearcoat system which delivered $800M revenue growth in 2012.'),\n    ([char]0x2022 + ' Developed new rheology test method with optimized shear profile, reducing error in measurement by 50%, for non' + [char]0xAD + '-Newtonian solventborne paint systems.'),\n    ([char]0x2022 + ' Technical lead for innovative spray process and paint technology conversion program which minimized assembly line downtime by 50% over conventional conversion.'),\n    ([char]0x2022 + ' Provided support to Arlington Assembly, GM''s most profitable manufacturing site, via new color development, formulation adjustments, and troubleshooting line issues (2011-2013). Currently provide support to Bowling Green Assembly, home of the Corvette.'),\n    ([char]0x2022 + ' Align formulas, manufacturing procedures, and product design specifications for manufacturing scale' + [char]0xAD + '-up.'),\n    ([char]0x2022 + ' Gravimetric measurement of binary vapor' + [char]0xAD + '-liquid equilibrium curves of ionic liquids with components in flue gas (e.g. CO2, CH4, H2O) as well as N2O.'),\n    ([char]0x2022 + ' Analyzed and calculated hysteresis, infinite dilution activity coefficients, Henry''s Law constants, and deconvolution of physical CO2 solubility from chemical CO2 reaction in amine' + [char]0xAD + '-functionalized ionic liquid systems.'),\n    ([char]0x2022 + ' Supervised design and construction of ionic liquid absorber/ stripper unit.'),\n    ([char]0x2022 + ' Proficient with both high and low pressure systems.'),\n    ([char]0x2022 + ' Facilitated formulation of controlled' + [char]0xAD + '-release coating.'),\n    ([char]0x2022 + ' Authored process flow diagrams. Collaborated on design of specialized fluidized bed dryer. Designated process instrumentation for final scale' + [char]0xAD + 'up.'),\n    ([char]0x2022 + ' Supervised 3' + [char]0xAD + '-4 production workers per shift in pilot plant operations.'),\n    ([char]0x2022 + ' Designed and formulated a novel filtration system capable of 100% toxin removal from contaminated liquids.'),\n    ([char]0x2022 + ' Bottled beer (KY Ale, KY Light, KY Bourbon Barrel Ale) at company microbrewery.')\n)\n\n$break = [char]0x0B\n\nforeach ($bulletText in $bulletsToRemove) {\n    $searchText = $bulletText + $break\n\n    $findRange = $d.Content\n    $find = $findRange.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.MatchWildcards = $false\n    $find.MatchCase = $true\n    $found = $find.Execute()\n\n    if (-not $found) {\n        throw \"Bullet text not found: $($bulletText.Substring(0, [Math]::Min(50, $bulletText.Length)))\"\n    }\n\n    $findRange.Delete()\n}\n"}
